$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.800.71'
$ws.Range("E2").Value = '  -0.08%  '

# Row 3
$ws.Range("D3").Value = '2.464.42'
$ws.Range("E3").Value = '  +0.80%  '

# Row 5
$ws.Range("D5").Value = '''574.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.36%  '

# Row 6
$ws.Range("D6").Value = '''147.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.27%  '

# Row 7
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").Value = '''0.534'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.17%  '

# Row 9
$ws.Range("D9").Value = '2.464.08'
$ws.Range("E9").Value = '  +0.87%  '

# Row 10
$ws.Range("E10").Value = '  +0.52%  '

# Row 11
$ws.Range("E11").Value = '  -0.49%  '

# Row 12
$ws.Range("D12").Value = '''5.27'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.04%  '

# Row 13
$ws.Range("E13").Value = '  +0.61%  '

# Row 14
$ws.Range("D14").Value = '''29.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.80%  '

# Row 15
$ws.Range("D15").Value = '''0.0000177'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.28%  '

# Row 16
$ws.Range("D16").Value = '2.912.00'

# Row 17
$ws.Range("D17").Value = '62.825.25'
$ws.Range("E17").Value = '  +0.06%  '

# Row 18
$ws.Range("D18").Value = '2.468.27'
$ws.Range("E18").Value = '  +0.94%  '

# Row 19
$ws.Range("E19").Value = '  +0.01%  '

# Row 20
$ws.Range("D20").Value = '''11.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.18%  '

# Row 21
$ws.Range("D21").Value = '''326.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.40%  '

# Row 22
$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").Value = '''4.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.07%  '

# Row 23
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").Value = '''2.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.92%  '

# Row 24
$ws.Range("E24").Value = '  -0.05%  '

# Row 25
$ws.Range("D25").Value = '''10.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +18.38%  '

# Row 26
$ws.Range("D26").Value = '''65.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.22%  '

# Row 27
$ws.Range("D27").Value = '''643.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.86%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0980'
$ws.Range("E28").Value = '  -1.28%  '

# Row 29
$ws.Range("D29").Value = '2.590.69'
$ws.Range("E29").Value = '  +0.95%  '

# Row 30
$ws.Range("D30").Value = '''0.997'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -15.30%  '

# Row 31
$ws.Range("E31").Value = '  -0.19%  '

# Row 32
$ws.Range("D32").Value = '''7.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.73%  '

# Row 33
$ws.Range("E33").Value = '  -1.59%  '

# Row 34
$ws.Range("E34").Value = '  -2.40%  '

# Row 35
$ws.Range("E35").Value = '  -0.02%  '

# Row 36
$ws.Range("D36").Value = '''1.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.29%  '

# Row 37
$ws.Range("E37").Value = '  -0.26%  '

# Row 38
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Value = '''0.369'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.38%  '

# Row 39
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Value = '''2.80'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.93%  '

# Row 40
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '''151.37'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.66%  '

# Row 41
$ws.Range("B41").Value = 'EthereumClassic'
$ws.Range("C41").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D41").Value = '''18.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.44%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").Value = '''5.38'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.56%  '

# Row 43
$ws.Range("E43").Value = '  -1.00%  '

# Row 44
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '''0.999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₆0307'
$ws.Range("E45").Value = '  -34.09%  '

# Row 46
$ws.Range("D46").Value = '''152.12'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.71%  '

# Row 47
$ws.Range("D47").Value = '''15.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.05%  '

# Row 48
$ws.Range("D48").Value = '''3.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.51%  '

# Row 49
$ws.Range("D49").Value = '''20.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.41%  '

# Row 50
$ws.Range("D50").Value = '''0.607'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.24%  '

# Row 51
$ws.Range("D51").Value = '''0.0511'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.97%  '
